# add-students.xlsx: split the multi-list "Faculty Advisor" batch-mapping
# helper columns (J:M on Sheet1) out onto their own "Sheet2", fix the
# mislabeled DSBS header, and add a brand-new "Career Option" lookup list
# (Superset Enrolled / Higher Studies / Entrepreneur) wired up as a second
# data-validation dropdown (column G).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Create Sheet2 right after Sheet1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# ---------------------------------------------------------------------
# 2. Move the CINTEL/DSBS/CTECH/NWC batch-mapping lists from Sheet1!J:M
#    onto Sheet2!A:D (only the cells that actually hold data).
# ---------------------------------------------------------------------
$srcCols = @("J", "K", "L", "M")
$dstCols = @("A", "B", "C", "D")
for ($i = 0; $i -lt $srcCols.Length; $i++) {
    for ($r = 1; $r -le 10; $r++) {
        $srcAddr = "$($srcCols[$i])$r"
        $dstAddr = "$($dstCols[$i])$r"
        $v = $ws1.Range($srcAddr).Value2
        if ($v -ne $null -and $v -ne "") {
            $ws2.Range($dstAddr).Value2 = $v
        }
    }
}

# The old K column header was a copy/paste bug (duplicated "CINTEL" instead
# of "DSBS") - fix it now that it lives in Sheet2!B1.
$ws2.Range("B1").Value2 = "DSBS"

# Remove the now-empty helper columns from Sheet1.
$ws1.Range("J1:M10").Clear()

# ---------------------------------------------------------------------
# 3. New "Career Option" lookup list on Sheet2!E1:E4
# ---------------------------------------------------------------------
$ws2.Range("E1").Value2 = "Career Option"
$ws2.Range("E2").Value2 = "Superset Enrolled"
$ws2.Range("E3").Value2 = "Higher Studies"
$ws2.Range("E4").Value2 = "Entrepreneur"

$fontE2 = $ws2.Range("E2").Font
$fontE2.Name = "Arial"
$fontE2.Size = 10
$fontE2.Color = 0

$fontE4 = $ws2.Range("E4").Font
$fontE4.Name = "Arial"
$fontE4.Size = 10
$fontE4.Color = 4473924

# ---------------------------------------------------------------------
# 4. Sheet2 cosmetics - column widths
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 29.5703125
$ws2.Columns.Item(2).ColumnWidth = 82.28515625
$ws2.Columns.Item(3).ColumnWidth = 37.140625
$ws2.Columns.Item(4).ColumnWidth = 74.42578125
$ws2.Columns.Item(5).ColumnWidth = 38.7109375

# ---------------------------------------------------------------------
# 5. Defined names now point at Sheet2 instead of Sheet1
# ---------------------------------------------------------------------
$wb.Names.Item("CINTEL").RefersTo = "=Sheet2!`$A`$2:`$A`$8"
$wb.Names.Item("CTECH").RefersTo = "=Sheet2!`$C`$2"
$wb.Names.Item("DSBS").RefersTo = "=Sheet2!`$B`$2:`$B`$7"
$wb.Names.Item("NWC").RefersTo = "=Sheet2!`$D`$2:`$D`$10"

# ---------------------------------------------------------------------
# 6. Data validations on Sheet1: repoint the Department (D) list at
#    Sheet2, and add the new Career Option (G) list.
# ---------------------------------------------------------------------
$ws1.Range("D1:D1048576").Validation.Delete()
$ws1.Range("D1:D1048576").Validation.Add(3, 1, 1, "=Sheet2!`$A`$1:`$D`$1")

$ws1.Range("G1:G1048576").Validation.Add(3, 1, 1, "=Sheet2!`$E`$2:`$E`$4")

# ---------------------------------------------------------------------
# 7. Selections: Sheet2 remembers its own scroll/selection state, but
#    Sheet1 stays the active tab with its selection moved to G5.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("F8").Select()

$ws1.Activate()
$ws1.Range("G5").Select()
